# Add two new columns "I0" (I) and "IF" (J) to the sheet, mirroring the
# style of the existing header row and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: same style as the other header cells (B1:H1 use style index 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J66 (column I then column J) per row.
$values = @(
    @(8, 8),    # row 2
    @(8, 8),    # row 3
    @(8, 8),    # row 4
    @(8, 8),    # row 5
    @(7, 8),    # row 6
    @(8, 8),    # row 7
    @(8, 8),    # row 8
    @(8, 8),    # row 9
    @(8, 8),    # row 10
    @(8, 8),    # row 11
    @(8, 8),    # row 12
    @(8, 8),    # row 13
    @(9, 9),    # row 14
    @(8, 8),    # row 15
    @(8, 8),    # row 16
    @(7, 8),    # row 17
    @(9, 9),    # row 18
    @(8, 8),    # row 19
    @(8, 8),    # row 20
    @(6, 6),    # row 21
    @(8, 8),    # row 22
    @(7, 7),    # row 23
    @(9, 9),    # row 24
    @(7, 7),    # row 25
    @(8, 8),    # row 26
    @(11, 11),  # row 27
    @(9, 9),    # row 28
    @(9, 9),    # row 29
    @(8, 9),    # row 30
    @(7, 7),    # row 31
    @(8, 8),    # row 32
    @(8, 8),    # row 33
    @(8, 9),    # row 34
    @(9, 9),    # row 35
    @(7, 7),    # row 36
    @(9, 9),    # row 37
    @(8, 8),    # row 38
    @(9, 9),    # row 39
    @(9, 9),    # row 40
    @(8, 8),    # row 41
    @(10, 10),  # row 42
    @(9, 9),    # row 43
    @(10, 10),  # row 44
    @(8, 8),    # row 45
    @(7, 7),    # row 46
    @(6, 6),    # row 47
    @(9, 9),    # row 48
    @(7, 7),    # row 49
    @(7, 8),    # row 50
    @(5, 6),    # row 51
    @(8, 8),    # row 52
    @(8, 8),    # row 53
    @(9, 9),    # row 54
    @(6, 6),    # row 55
    @(9, 9),    # row 56
    @(9, 9),    # row 57
    @(8, 8),    # row 58
    @(8, 8),    # row 59
    @(8, 8),    # row 60
    @(7, 7),    # row 61
    @(8, 8),    # row 62
    @(7, 7),    # row 63
    @(8, 8),    # row 64
    @(6, 6),    # row 65
    @(7, 7)     # row 66
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 9).Value = $values[$i][0]
    $ws.Cells.Item($r, 10).Value = $values[$i][1]
}
